# Incorporated old tests into new scheme (#225)
#
# This script reproduces, via the Excel object model, the changes that were
# made to zigar-compiler/test/integration/test-matrix.xlsx:
#   - "Type handling"     : just move the active-cell selection (E32 -> J14)
#   - "Error handling"    : drop the trailing blank row (dimension A1:B4 -> A1:B3)
#   - "Function calling"  : widen column A and append 11 new scenario rows
#                           (dimension A1:B4 -> A1:B14)
#   - "Memory allocation" : relabel the existing scenario row and fill in the
#                           previously-blank trailing row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Error handling": the trailing, empty row 4 is removed.
# ---------------------------------------------------------------------------
$wsError = $wb.Worksheets.Item("Error handling")
$wsError.Rows.Item(4).Delete()
$wsError.Range("A4").Select()

# ---------------------------------------------------------------------------
# Sheet "Function calling": column A grows wider and 11 new scenario rows
# (rows 4-14) are filled in below the existing "Attach getters/setters" row.
# ---------------------------------------------------------------------------
$wsFunc = $wb.Worksheets.Item("Function calling")

# width="29.1" -> width="35.06" (nearest value obtainable through the
# pixel-quantized ColumnWidth setter)
$wsFunc.Columns.Item(1).ColumnWidth = 34.15

# Row 3's scenario slot is relabelled from "Attach getters/setters" to
# "Throw error"; "Attach getters/setters" reappears later, as row 14.
$wsFunc.Range("A3").Value = "Throw error"

# Row 4 already exists (with the right cell styles) in the workbook, so its
# value can just be assigned directly.
$wsFunc.Range("A4").Value = "Return slice"
$wsFunc.Range("B4").Value = "Y"

# Rows 5-14 are brand new, so copy the formatting from row 4 first and then
# fill in the labels.
$newScenarios = @(
    "Return slice of slices",
    "Print slice of slices",
    "Accept typed array",
    "Return bool vector",
    "Handle misaligned pointer",
    "Handled misaligned aliased pointer",
    "Allocate slice of structs",
    "Clear pointers",
    "Clear pointer array",
    "Attach getters/setters"
)

$wsFunc.Range("A4:B4").Copy()
for ($i = 0; $i -lt $newScenarios.Length; $i++) {
    $row = 5 + $i
    $wsFunc.Range("A" + $row + ":B" + $row).PasteSpecial(-4122)
    $wsFunc.Range("A" + $row).Value = $newScenarios[$i]
    $wsFunc.Range("B" + $row).Value = "Y"
}

$wsFunc.Range("E34").Select()

# ---------------------------------------------------------------------------
# Sheet "Memory allocation": the existing row 3 keeps its text ("Create
# internal slice") but row 4, formerly blank, now holds a new scenario.
# ---------------------------------------------------------------------------
$wsMem = $wb.Worksheets.Item("Memory allocation")
$wsMem.Range("A3").Value = "Create internal slice"
$wsMem.Range("A4").Value = "Allocate memory for string"
$wsMem.Range("B4").Value = "Y"
$wsMem.Range("A5").Select()

# ---------------------------------------------------------------------------
# Sheet "Type handling": only the remembered selection changes. This is done
# last so that "Type handling" ends up the active/selected sheet again, just
# like in the original workbook.
# ---------------------------------------------------------------------------
$wsType = $wb.Worksheets.Item("Type handling")
$wsType.Activate()
$wsType.Range("J14").Select()
